$d = $word.ActiveDocument

# --- 1. Remove the trailing two Title paragraphs ("Restrições..." and "Esclarecimentos") ---
$d.Paragraphs(4).Range.Delete()
$d.Paragraphs(3).Range.Delete()

# --- 2. Remove the long confidentiality body paragraph (paragraph 2) ---
$d.Paragraphs(2).Range.Delete()

# --- 3. Turn the remaining paragraph into the new "Section 2: Data" Heading2 ---
$p1 = $d.Paragraphs(1)
$headingXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Section 2: Data</w:t></w:r></w:p>'
[void]$p1.Range.InsertXML($headingXml)

# --- 4. Insert the new 2-column support table right after that paragraph ---
$insertPos = $d.Paragraphs(1).Range.End
$insertRange = $d.Range($insertPos, $insertPos)

$tableXml = '<w:tbl xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:tblPr><w:tblStyle w:val="TableGrid"/><w:tblW w:w="0" w:type="auto"/><w:tblLayout w:type="fixed"/><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="4320"/><w:gridCol w:w="4320"/></w:tblGrid><w:tr><w:tc><w:tcPr><w:tcW w:w="4000" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>&#8226;</w:t><w:tab/><w:t>Acesso pelo telefone a especialistas</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4000" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>&#8226;</w:t><w:tab/><w:t>Alertas preventivos HPE InfoSight</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="4000" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>&#8226;</w:t><w:tab/><w:t>Bate-papo online com especialistas</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4000" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>&#8226;</w:t><w:tab/><w:t>Registro de incidentes automatizado</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="4000" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>&#8226;</w:t><w:tab/><w:t>Respostas ao fórum dadas por especialistas</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4000" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>&#8226;</w:t><w:tab/><w:t>Biblioteca de dicas técnicas</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="4000" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>&#8226;</w:t><w:tab/><w:t>Orientação técnica geral</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4000" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>&#8226;</w:t><w:tab/><w:t>Acesso a informações e serviços de suporte eletrônico</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="4000" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>&#8226;</w:t><w:tab/><w:t>Assistência HPE InfoSight</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4000" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>&#8226;</w:t><w:tab/><w:t>Gerenciamento de interrupções (apenas no nível de serviço de Crítico)</w:t></w:r></w:p></w:tc></w:tr></w:tbl>'

[void]$insertRange.InsertXML($tableXml)

Write-Host "Paragraphs: " $d.Paragraphs.Count
Write-Host "Tables: " $d.Tables.Count
